# Refresh cryptos list with latest scraped price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.209.76'
$ws.Range("E2").Value = '  -1.58%  '
$ws.Range("D3").Value = '2.245.34'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'246.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("E6").Value = '  -1.43%  '
$ws.Range("D7").Value = "'74.31"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.91%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -3.92%  '
$ws.Range("D10").Value = "'42.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.18%  '
$ws.Range("D11").Value = "'0.0948"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.64%  '
$ws.Range("D12").Value = "'7.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("E14").Value = '  -4.86%  '
$ws.Range("E15").Value = '  -2.08%  '
$ws.Range("D16").Value = '2.235.00'
$ws.Range("E16").Value = '  -1.96%  '
$ws.Range("D17").Value = '42.103.31'
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("D18").Value = '0.0₃0994'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = "'72.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").Value = "'6.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.53%  '
$ws.Range("D21").Value = "'2.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.98%  '
$ws.Range("D22").Value = "'230.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("D23").Value = "'9.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +36.72%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("E26").Value = '  -4.38%  '
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("D28").Value = "'2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.55%  '
$ws.Range("D29").Value = "'171.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.79%  '
$ws.Range("D30").Value = "'20.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  -4.29%  '
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("D33").Value = "'30.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.80%  '
$ws.Range("E34").Value = '  +10.35%  '
$ws.Range("E35").Value = '  -0.82%  '
$ws.Range("E36").Value = '  -2.78%  '
$ws.Range("D37").Value = "'0.0314"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.30%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("E39").Value = '  -4.03%  '
$ws.Range("E40").Value = '  -1.73%  '
$ws.Range("D41").Value = "'62.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.25%  '
$ws.Range("E42").Value = '  -2.35%  '
$ws.Range("D43").Value = "'106.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = '  +1.68%  '
$ws.Range("E45").Value = '  -2.89%  '
$ws.Range("D46").Value = "'0.996"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").Value = "'1.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = "'1.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.37%  '
$ws.Range("D49").Value = "'2.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.85%  '
$ws.Range("E50").Value = '  -7.66%  '
